$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ExpenseRequest")

# Update the Event name value used in the test data (shared string update)
$ws.Range("D2").Value = "PFG Golf Event"
$ws.Range("D3").Value = "PFG Golf Event"

# Move the active cell selection to D13 as recorded when the sheet was last saved
$ws.Activate()
$ws.Range("D13").Select()
